# Add new treasury/fed columns (F:J) to the legal-dates sheet and populate
# them for rows 2-21, matching the upstream data arrangement fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (F1:J1): clone E1's header style, then set the new labels ---
$ws.Range("E1").Copy()
$ws.Range("F1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F1").Value = "treasury_open"
$ws.Range("G1").Value = "treasury_close"
$ws.Range("H1").Value = "treasury_delta"
$ws.Range("I1").Value = "fed_maturities"
$ws.Range("J1").Value = "fed_investments"

# --- Data rows 2-17: treasury_open, treasury_close, treasury_delta, fed_maturities, fed_investments ---
$treasuryFedData = @{
  2  = @(1624404,  1642285,  17881,    10875237800, 6001000000)
  3  = @(1642285,  1641080,  -1205,    0,           1201000000)
  4  = @(0,        0,        0,        0,           0)
  5  = @(0,        0,        0,        0,           0)
  6  = @(1641080,  1653555,  12475,    0,           0)
  7  = @(1653555,  1627709,  -25846,   6657527900,  1732000000)
  8  = @(1627709,  1740455,  112746,   7457547200,  12801000000)
  9  = @(1740455,  1807305,  66850,    18724230800, 3601000000)
  10 = @(1807305,  1812801,  5496,     0,           1199000000)
  11 = @(0,        0,        0,        0,           0)
  12 = @(0,        0,        0,        0,           0)
  13 = @(1812801,  1817270,  4469,     0,           0)
  14 = @(1817270,  1794448,  -22822,   5332930800,  1734000000)
  15 = @(1794448,  1777180,  -17268,   0,           6001000000)
  16 = @(1777180,  1821824,  44644,    16277175800, 2401000000)
  17 = @(1821824,  1825498,  3674,     0,           1732000000)
}

foreach ($row in 2..17) {
  $vals = $treasuryFedData[$row]
  $ws.Range("F$row").Value = $vals[0]
  $ws.Range("G$row").Value = $vals[1]
  $ws.Range("H$row").Value = $vals[2]
  $ws.Range("I$row").Value = $vals[3]
  $ws.Range("J$row").Value = $vals[4]
}

# --- Rows 18-21: treasury columns (F:G:H) are blank text cells, fed columns (I:J) stay numeric ---
$fedOnlyData = @{
  18 = @(0,          0)
  19 = @(0,          0)
  20 = @(0,          8801000000)
  21 = @(6398553400, 12825000000)
}

foreach ($row in 18..21) {
  foreach ($col in "F", "G", "H") {
    # Assign via a leading apostrophe so the cell stays text-typed but empty,
    # then reset the style back to Normal (no quote-prefix formatting left behind).
    $ws.Range("$col$row").Value = "'"
    $ws.Range("$col$row").Style = "Normal"
  }
  $vals = $fedOnlyData[$row]
  $ws.Range("I$row").Value = $vals[0]
  $ws.Range("J$row").Value = $vals[1]
}

Write-Host "done"
